# Generate Report for Handback
#
# - Marks the zh-cn / de-de handoffs as handed back (Status text + Latest
#   Handback DateTime).
# - Fills in the "Latest Target File" / "Latest Handback File" columns that
#   were previously blank, with hyperlinks pointing at the same targets
#   used elsewhere on each row.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1. Status column: "Ready for handoff" -> "Handed back: in sync with en-US"
# ---------------------------------------------------------------------
$wsOverview.Range("B2").Value = $statusHandedBack
$wsOverview.Range("C2").Value = $statusHandedBack
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

$wsZhCn.Range("C2").Value = $statusHandedBack
$wsZhCn.Range("C3").Value = $statusHandedBack

$wsDeDe.Range("C2").Value = $statusHandedBack
$wsDeDe.Range("C3").Value = $statusHandedBack

# ---------------------------------------------------------------------
# 2. Latest Handback DateTime (column H)
# ---------------------------------------------------------------------
$wsZhCn.Range("H2").Value = "2016-03-21 08:57:12"
$wsZhCn.Range("H3").Value = "2016-03-21 08:57:12"

$wsDeDe.Range("H2").Value = "2016-03-21 08:57:19"
$wsDeDe.Range("H3").Value = "2016-03-21 08:57:19"

# ---------------------------------------------------------------------
# 3. Latest Target File (F) / Latest Handback File (G) + hyperlinks
# ---------------------------------------------------------------------
$mdName   = "227dfc44-15b6-42d6-a696-ae27c012d273.md"
$mdUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/80525a9c48e0f1002cfddc3c207c02c24b6d4211/e2e/227dfc44-15b6-42d6-a696-ae27c012d273.md"

$zhCnXlfName = "227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.zh-cn.xlf"
$zhCnXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5dd169b752834a2f3d6b9be2527ab9ae06400fc5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.zh-cn.xlf"

$deDeXlfName = "227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.de-de.xlf"
$deDeXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6a100a0e16111973fd8f31dab16272974c9c453e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.de-de.xlf"

# Hyperlink style colour (matches the workbook's existing "HyperLink" cell
# style: underline, font colour #6495ED). Excel's Font.Color is a BGR-packed
# decimal, hence the byte-swapped literal below.
$hyperlinkColor = 15570276

function Add-HandbackHyperlink($ws, $cellRef, $url, $displayText) {
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, "", "", $displayText) | Out-Null
    $ws.Range($cellRef).Font.Underline = $true
    $ws.Range($cellRef).Font.Color = $hyperlinkColor
}

# zh-cn: row 2
Add-HandbackHyperlink $wsZhCn "F2" $mdUrl $mdName
Add-HandbackHyperlink $wsZhCn "G2" $zhCnXlfUrl $zhCnXlfName
# zh-cn: row 3
Add-HandbackHyperlink $wsZhCn "F3" $mdUrl $mdName
Add-HandbackHyperlink $wsZhCn "G3" $zhCnXlfUrl $zhCnXlfName

# de-de: row 2
Add-HandbackHyperlink $wsDeDe "F2" $mdUrl $mdName
Add-HandbackHyperlink $wsDeDe "G2" $deDeXlfUrl $deDeXlfName
# de-de: row 3
Add-HandbackHyperlink $wsDeDe "F3" $mdUrl $mdName
Add-HandbackHyperlink $wsDeDe "G3" $deDeXlfUrl $deDeXlfName

Write-Output "Handback report generated."
